$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new worker-detail row (row 18), pushing the signature block
#     (old rows 22/23) down to rows 23/24, exactly as Excel does when a full
#     row is inserted above them. ---
$ws.Cells.Item(18, 1).EntireRow.Insert()

# Row 17 (the former last/bottom data row, heavier "closing" border style)
# still holds its original content at this point; duplicate it into the
# freshly inserted row 18 (value + full formatting, via Copy(Destination)
# which -- unlike Copy()+Paste -- reliably carries borders/fill/font too).
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Row 16 (the first data row style) gets duplicated into row 17, which
# overwrites row 17's old content/style with row 16's lighter style -- this
# mirrors how the new period row was actually inserted above the old
# bottom row in the source workbook.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# --- Now set the three "Periodo Mora" values for the three worker rows. ---
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2507"
$ws.Range("E18").Value = "2508"

# --- Update the summary figures: VALOR MORA total and Cant. Periodos count. ---
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
